$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values look numeric (e.g. "217.30"); setting .Value directly
# would let Excel auto-coerce them to actual numbers and drop formatting like
# trailing zeros or the thousands-style dots used by coinranking.com exports.
# Force text entry via NumberFormat="@", then restore the cell style to Normal
# (matching the original, un-styled inlineStr cells) so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.842.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.648.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.252"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0629"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.871.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.638.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.530"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.844.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.78%  "
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.291.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("E37").Value = "  -5.81%  "
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.829"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.797.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0979"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "
